$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.159.81"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.264.48"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'306.73"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'96.77"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("D7").Value = "'0.525"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "'35.10"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'6.88"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "2.615.96"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'14.62"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "2.271.41"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "42.036.72"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'12.39"
$ws.Range("E19").Value = "  -4.27%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "'6.00"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "'68.30"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'237.53"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "'2.57"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'23.57"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'37.19"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.12"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.48"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "'160.44"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'5.21"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'3.19"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").Value = "'0.0736"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "'17.23"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "'0.104"
$ws.Range("E38").Value = "  -3.27%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'3.98"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("D42").Value = "'2.31"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").Value = "1.957.25"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'18.98"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").Value = "'0.0281"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'10.00"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'53.22"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'71.61"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'91.79"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -1.44%  "
